$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old "errors" rows (3 and 4) entirely - they are not part of the new data
$null = $ws.Range("A3:B4").ClearContents()

# Row 1: search term + matched product link (logged search result #1)
$ws.Range("A1").Value = "samsung galaxy fold negro"
$ws.Range("B1").Value = "https://www.fnac.es/Samsung-Galaxy-Fold-7-3-512GB-Negro-Telefono-movil-Smartphone/a7132771#omnsearchpos=11"

# Row 2: search term + matched product link (logged search result #2)
$ws.Range("A2").Value = "iphone 12 pro max oro"
$ws.Range("B2").Value = "https://www.backmarket.es/iphone-12-pro-max-128-gb-oro-libre-segunda-mano/413839.html#l=10"

# Row 15/16: the same two log entries re-appended further down (e.g. a logging
# function that appends each run's results to the bottom of the sheet)
$ws.Range("A15").Value = "samsung galaxy fold negro"
$ws.Range("B15").Value = "https://www.fnac.es/Samsung-Galaxy-Fold-7-3-512GB-Negro-Telefono-movil-Smartphone/a7132771#omnsearchpos=11"

$ws.Range("A16").Value = "iphone 12 pro max oro"
$ws.Range("B16").Value = "https://www.backmarket.es/iphone-12-pro-max-128-gb-oro-libre-segunda-mano/413839.html#l=10"

# Mimic "Select All" (Ctrl+A) being the active selection when the file was saved
$null = $ws.Cells.Select()
